$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 203, shifting rows 203:258 down to 204:259
$ws.Rows(203).Insert()

# Populate the newly inserted row 203 with its values
$ws.Range("A203").Value = 4
$ws.Range("B203").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C203").Value = "Los Lagos"
$ws.Range("D203").Value = 44841
$ws.Range("E203").Value = 10
$ws.Range("F203").Value = "Fruta"
$ws.Range("G203").Value = 100108
$ws.Range("H203").Value = "Tropicales y subtropicales"
$ws.Range("I203").Value = 100108002
$ws.Range("J203").Value = "Mango"
$ws.Range("K203").Value = "Sin especificar"
$ws.Range("L203").Value = "Primera"
$ws.Range("M203").Value = 200
$ws.Range("N203").Value = 9000
$ws.Range("O203").Value = 10000
$ws.Range("P203").Value = 9500
$ws.Range("Q203").Value = "$/bandeja 4 kilos"
$ws.Range("R203").Value = "Brasil"
$ws.Range("S203").Value = 2375
$ws.Range("T203").Value = 4
